# Daily attendance processing - 2025-12-08 11:27:23
#
# Column G ("Recorded By") holds a comma-separated list of the
# users/systems that recorded/touched each attendance session row.
# These lists need to be normalised into a deterministic (ordinal,
# case-sensitive, ASCII) sorted order, e.g.
#   "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
#   "system, System, backup@backdoor.com" -> "System, backup@backdoor.com, system"
#
# PowerShell's own Sort-Object / -ceq / -clt operators in this host are
# culture-aware (case-insensitive), so we roll a tiny ordinal
# string-compare + insertion sort by hand using character codes.

function OrdCompare($ocX, $ocY) {
    $ocLx = $ocX.Length
    $ocLy = $ocY.Length
    $ocN = $ocLx
    if ($ocLy -lt $ocN) { $ocN = $ocLy }
    for ($ocI = 0; $ocI -lt $ocN; $ocI++) {
        $ocCx = [int][char]$ocX[$ocI]
        $ocCy = [int][char]$ocY[$ocI]
        if ($ocCx -lt $ocCy) { return -1 }
        if ($ocCx -gt $ocCy) { return 1 }
    }
    if ($ocLx -lt $ocLy) { return -1 }
    if ($ocLx -gt $ocLy) { return 1 }
    return 0
}

function OrdSort($osArr) {
    $osN = $osArr.Length
    for ($osI = 1; $osI -lt $osN; $osI++) {
        $osKey = $osArr[$osI]
        $osJ = $osI - 1
        while ($osJ -ge 0 -and (OrdCompare $osArr[$osJ] $osKey) -gt 0) {
            $osArr[$osJ + 1] = $osArr[$osJ]
            $osJ = $osJ - 1
        }
        $osArr[$osJ + 1] = $osKey
    }
    return $osArr
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $firstRow + $usedRange.Rows.Count - 1

# Row 1 is the header ("Recorded By"); data starts at row 2.
$dataFirstRow = $firstRow
if ($dataFirstRow -lt 2) { $dataFirstRow = 2 }

for ($rowNum = $dataFirstRow; $rowNum -le $lastRow; $rowNum++) {
    $cell = $ws.Cells.Item($rowNum, 7)
    $originalText = $cell.Value2

    if ($originalText -ne $null -and $originalText -ne "") {
        $parts = $originalText.Split(",")
        $trimmedList = @()
        foreach ($piece in $parts) {
            $trimmedList += $piece.Trim()
        }

        $sortedList = OrdSort $trimmedList
        $newText = [string]::Join(", ", $sortedList)

        if ($newText -cne $originalText) {
            $cell.Value = $newText
        }
    }
}
